$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells AD1:AF1 ("Wins", "Losses", "Ties"), styled like the
#     rest of the header row (bold, thin border, centered/top aligned) ---
$headers = @{ "AD1" = "Wins"; "AE1" = "Losses"; "AF1" = "Ties" }
foreach ($addr in $headers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $headers[$addr]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# --- Season record (Wins/Losses/Ties) for every data row, 2 through 41 ---
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 96   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 66   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
